$d = $word.ActiveDocument

# The last paragraph in the document is the empty paragraph that sits
# right after the "Algoritmo recursivo (força bruta)" heading. We fill
# it in with the explanatory text, using subscripts for the "i"/"j"
# index variables (xi, xj) exactly like the target revision.

$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)
$range = $target.Range
$range.Collapse(1)

$fullText = "O algoritmo recursivo de força bruta tem como caso base se I = j = 0, nesse caso o algoritmo retorna 0. Caso contrário, ele tem dois caso não base, o primeiro, quando xi " + [char]0x2260 + " xj ele retorna max(LCS(i-1,j),LCS(i,j-1)), se xi=xj ele retorna LCS(i-1,j-1)+1."

$range.InsertAfter($fullText)

$paraRange = $d.Paragraphs.Item($count).Range
$base = $paraRange.Start

function Set-Sub($from, $to) {
    $r = $d.Range($base + $from, $base + $to)
    $r.Font.Subscript = $true
}

function Set-Baseline($from, $to) {
    $r = $d.Range($base + $from, $base + $to)
    $r.Font.Subscript = $false
}

# offsets computed from the text above:
# 0-168   : normal text ending in "...quando x"
# 168-170 : "i " (subscript)
# 170-173 : "≠ x" (baseline)
# 173-174 : "j" (subscript)
# 174-219 : " ele retorna max(LCS(i-1,j),LCS(i,j-1)), se x" (baseline)
# 219-220 : "i" (subscript)
# 220-222 : "=x" (baseline)
# 222-223 : "j" (subscript)
# 223-251 : " ele retorna LCS(i-1,j-1)+1." (baseline)

function Check-Segment($from, $to, $expected) {
    $r = $d.Range($base + $from, $base + $to)
    if ($r.Text -ne $expected) {
        Write-Output "MISMATCH at $from-$to : got [$($r.Text)] expected [$expected]"
    }
}

$neq = [string]([char]0x2260) + " x"

Check-Segment 168 170 "i "
Check-Segment 170 173 $neq
Check-Segment 173 174 "j"
Check-Segment 174 219 " ele retorna max(LCS(i-1,j),LCS(i,j-1)), se x"
Check-Segment 219 220 "i"
Check-Segment 220 222 "=x"
Check-Segment 222 223 "j"
Check-Segment 223 251 " ele retorna LCS(i-1,j-1)+1."

Set-Sub 168 170
Set-Baseline 170 173
Set-Sub 173 174
Set-Baseline 174 219
Set-Sub 219 220
Set-Baseline 220 222
Set-Sub 222 223
Set-Baseline 223 251

Write-Output "done: $($d.Paragraphs.Item($count).Range.Text)"
